# Update crypto price/volume data to match the Sat Jan 27 11:54:19 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.771.45"
$ws.Range("E2").Value = "'  +1.35%  "
$ws.Range("D3").Value = "'2.268.50"
$ws.Range("E3").Value = "'  +0.66%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'304.46"
$ws.Range("E5").Value = "'  +0.58%  "
$ws.Range("D6").Value = "'91.93"
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("E7").Value = "'  +1.65%  "
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E9").Value = "'  -0.30%  "
$ws.Range("D10").Value = "'32.32"
$ws.Range("E10").Value = "'  +0.37%  "
$ws.Range("D11").Value = "'53.32"
$ws.Range("E11").Value = "'  +1.26%  "
$ws.Range("E12").Value = "'  +0.29%  "
$ws.Range("E13").Value = "'  -0.09%  "
$ws.Range("E14").Value = "'  +0.90%  "
$ws.Range("D15").Value = "'2.618.21"
$ws.Range("E15").Value = "'  +0.59%  "
$ws.Range("E16").Value = "'  +0.76%  "
$ws.Range("D17").Value = "'2.262.50"
$ws.Range("E17").Value = "'  +1.50%  "
$ws.Range("D18").Value = "'0.765"
$ws.Range("E18").Value = "'  +1.85%  "
$ws.Range("D19").Value = "'41.690.65"
$ws.Range("E19").Value = "'  +1.37%  "
$ws.Range("E20").Value = "'  +5.55%  "
$ws.Range("E21").Value = "'  +0.15%  "
$ws.Range("D22").Value = "'5.93"
$ws.Range("E22").Value = "'  +1.05%  "
$ws.Range("D23").Value = "'67.07"
$ws.Range("E23").Value = "'  +0.46%  "
$ws.Range("D24").Value = "'239.96"
$ws.Range("E24").Value = "'  +0.21%  "
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "'  +0.83%  "
$ws.Range("E26").Value = "'  +0.03%  "
$ws.Range("E27").Value = "'  +2.21%  "
$ws.Range("D28").Value = "'23.96"
$ws.Range("E28").Value = "'  +0.09%  "
$ws.Range("D29").Value = "'9.53"
$ws.Range("E29").Value = "'  -0.12%  "
$ws.Range("D30").Value = "'2.06"
$ws.Range("E30").Value = "'  -4.93%  "
$ws.Range("D31").Value = "'34.79"
$ws.Range("E31").Value = "'  +3.84%  "
$ws.Range("D32").Value = "'160.32"
$ws.Range("E32").Value = "'  +0.97%  "
$ws.Range("E33").Value = "'  +2.64%  "
$ws.Range("E34").Value = "'  -0.18%  "
$ws.Range("E35").Value = "'  +1.51%  "
$ws.Range("E36").Value = "'  -1.40%  "
$ws.Range("D37").Value = "'16.80"
$ws.Range("E37").Value = "'  +2.53%  "
$ws.Range("E38").Value = "'  +0.54%  "
$ws.Range("E39").Value = "'  +1.29%  "
$ws.Range("E40").Value = "'  -0.61%  "
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "'  +0.61%  "
$ws.Range("E42").Value = "'  +0.14%  "
$ws.Range("D43").Value = "'2.026.31"
$ws.Range("E43").Value = "'  -2.74%  "
$ws.Range("D44").Value = "'19.30"
$ws.Range("E44").Value = "'  -3.80%  "
$ws.Range("D45").Value = "'10.41"
$ws.Range("E45").Value = "'  +0.20%  "
$ws.Range("E46").Value = "'  +0.79%  "
$ws.Range("E47").Value = "'  +12.77%  "
$ws.Range("E48").Value = "'  -1.58%  "
$ws.Range("B49").Value = "'Stacks"
$ws.Range("C49").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.52"
$ws.Range("E49").Value = "'  -0.68%  "
$ws.Range("B50").Value = "'BitcoinSV"
$ws.Range("C50").Value = "'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'72.59"
$ws.Range("E50").Value = "'  +3.63%  "
$ws.Range("D51").Value = "'1.15"
$ws.Range("E51").Value = "'  +0.39%  "
